# Commit: for panels.py add 1.5 for ultimate loads
# Updates the "RF" (reserve factor) columns (G, M, S) on rows 2-11 to reflect
# the new ultimate-load factor of 1.5 applied in panels.py.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.01
$ws.Range("M2").Value = 0.96
$ws.Range("S2").Value = 0.23

$ws.Range("G3").Value = 1.01
$ws.Range("M3").Value = 0.96
$ws.Range("S3").Value = 0.23

$ws.Range("G4").Value = 0.76
$ws.Range("M4").Value = 0.8100000000000001
$ws.Range("S4").Value = 0.22

$ws.Range("G5").Value = 0.76
$ws.Range("M5").Value = 0.8
$ws.Range("S5").Value = 0.22

$ws.Range("G6").Value = 0.76
$ws.Range("M6").Value = 0.8
$ws.Range("S6").Value = 0.22

$ws.Range("G7").Value = 0.65
$ws.Range("M7").Value = 0.74
$ws.Range("S7").Value = 0.08

$ws.Range("G8").Value = 0.65
$ws.Range("M8").Value = 0.74
$ws.Range("S8").Value = 0.08

$ws.Range("G9").Value = 0.76
$ws.Range("M9").Value = 0.8
$ws.Range("S9").Value = 0.22

$ws.Range("G10").Value = 0.76
$ws.Range("M10").Value = 0.8
$ws.Range("S10").Value = 0.22

$ws.Range("G11").Value = 0.76
$ws.Range("M11").Value = 0.8100000000000001
$ws.Range("S11").Value = 0.22
